$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" header on the "Weekly Quantity" sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Rename "Requested quantity" header on the "Monthly Trend" sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add a new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, bordered, centered) used by the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)
$wsForecast.Range("C1:D1").PasteSpecial(-4122)

# Match the date-time number format used for the "ds" column elsewhere
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A53").PasteSpecial(-4122)

# Forecast data rows
$data = @(
  @(44934.99999999999, 20, -2.001083464255852, 41.96366744955531),
  @(44941.99999999999, 20, -2.112570928467731, 42.09254713643674),
  @(44983.99999999999, 20, -3.282147916663185, 42.24414030731676),
  @(44990.99999999999, 20, -2.902805895993974, 42.41759671063445),
  @(44997.99999999999, 20, -2.042002740609933, 41.69574964538967),
  @(45004.99999999999, 20, -4.235735814107168, 42.02656732373559),
  @(45011.99999999999, 20, -2.264620392690848, 39.33128303614692),
  @(45018.99999999999, 20, -2.145640223112871, 42.27314951593049),
  @(45025.99999999999, 20, -3.611270092610765, 42.23106518668574),
  @(45039.99999999999, 20, -3.089362907673929, 40.19045518204013),
  @(45046.99999999999, 20, -1.266583693098193, 42.95876673060301),
  @(45053.99999999999, 19, -2.385294085100201, 42.16914752916026),
  @(45060.99999999999, 19, -2.254251291527396, 42.1315196626172),
  @(45067.99999999999, 19, -3.058169015468867, 41.41153487003594),
  @(45081.99999999999, 19, -3.534441180629383, 41.88467653650137),
  @(45088.99999999999, 19, -4.109763876675293, 41.66665781357219),
  @(45095.99999999999, 19, -3.511597366912603, 43.49980438871455),
  @(45102.99999999999, 19, -3.188688305520853, 41.52038194698891),
  @(45109.99999999999, 19, -3.544246417302229, 41.54621501799996),
  @(45116.99999999999, 19, -2.396730188976981, 41.32173261272515),
  @(45130.99999999999, 19, -3.331300059612547, 40.4430713162415),
  @(45151.99999999999, 19, -2.524428478156972, 41.99191117581515),
  @(45158.99999999999, 19, -4.326158616405931, 41.23295086797926),
  @(45186.99999999999, 18, -1.308524441105588, 42.06342542537123),
  @(45193.99999999999, 18, -2.728243237580706, 41.27232557866788),
  @(45200.99999999999, 18, -4.658265096277596, 41.89109501890778),
  @(45214.99999999999, 18, -4.328789426734062, 42.14271137590939),
  @(45221.99999999999, 18, -4.779903970789701, 40.33719225407743),
  @(45228.99999999999, 18, -4.074639311602143, 39.93532624021667),
  @(45235.99999999999, 18, -5.190385641307411, 41.01987262422956),
  @(45242.99999999999, 18, -3.510546647329922, 39.54281814454633),
  @(45249.99999999999, 18, -3.864464741188095, 40.48708271900549),
  @(45256.99999999999, 18, -5.225971706246503, 39.79409969105137),
  @(45263.99999999999, 18, -2.99653941279617, 40.41370352758963),
  @(45270.99999999999, 18, -4.572053819219883, 38.8296446354957),
  @(45277.99999999999, 18, -3.875548105249989, 40.72126439055647),
  @(45298.99999999999, 18, -2.394984604430462, 40.77689916758364),
  @(45305.99999999999, 17, -3.87380615422391, 39.02488823337442),
  @(45312.99999999999, 17, -6.330881602104323, 40.64618149951462),
  @(45319.99999999999, 17, -3.932065431696644, 39.9344703337876),
  @(45473.99999999999, 16, -6.889046017631285, 39.25516001879645),
  @(45487.99999999999, 16, -6.532220048094338, 39.21738200544963),
  @(45494.99999999999, 16, -7.165344900999286, 39.01019700591439),
  @(45501.99999999999, 16, -6.458806781291306, 40.3876557959607),
  @(45508.99999999999, 16, -6.796791392164205, 37.70709528177134),
  @(45515.99999999999, 16, -6.386553536273462, 39.00398611132584),
  @(45522.99999999999, 16, -7.023032179087068, 37.47406802633554),
  @(45529.99999999999, 16, -7.131361432447015, 39.11167892606952),
  @(45536.99999999999, 16, -5.73525563810245, 37.07065877135237),
  @(45543.99999999999, 16, -7.320946901800112, 36.92161248142936),
  @(45550.99999999999, 16, -6.195085520065345, 39.27765650385987),
  @(45557.99999999999, 15, -6.555332465829763, 38.72294947421716)
)

$r = 2
foreach ($row in $data) {
  $wsForecast.Cells.Item($r, 1).Value = $row[0]
  $wsForecast.Cells.Item($r, 2).Value = $row[1]
  $wsForecast.Cells.Item($r, 3).Value = $row[2]
  $wsForecast.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
